$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices formatted as text in the source data (values such as
# "69.154.50" or "1.00" are literal strings, not numbers/dates). Pre-format the
# cells we are about to rewrite as Text so Excel does not auto-convert them.
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D51").NumberFormat = "@"

# --- Column B (Coin) updates ---
$ws.Range("B20").Value = 'TRON'
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("B34").Value = 'Cosmos'
$ws.Range("B35").Value = 'OKB'
$ws.Range("B36").Value = 'Hedera'
$ws.Range("B39").Value = 'CoreDAO'
$ws.Range("B40").Value = 'PEPE'
$ws.Range("B41").Value = 'Stacks'
$ws.Range("B42").Value = 'TheGraph'
$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("B48").Value = 'ThetaToken'

# --- Column C (Link) updates ---
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("C48").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'

# --- Column D (Price) updates ---
$ws.Range("D2").Value = '69.154.50'
$ws.Range("D3").Value = '3.477.11'
$ws.Range("D5").Value = '574.53'
$ws.Range("D6").Value = '188.73'
$ws.Range("D7").Value = '3.467.29'
$ws.Range("D8").Value = '0.606'
$ws.Range("D9").Value = '1.00'
$ws.Range("D10").Value = '0.201'
$ws.Range("D11").Value = '0.616'
$ws.Range("D12").Value = '50.57'
$ws.Range("D13").Value = '0.0000281'
$ws.Range("D14").Value = '9.03'
$ws.Range("D15").Value = '4.037.81'
$ws.Range("D16").Value = '632.50'
$ws.Range("D17").Value = '69.006.71'
$ws.Range("D18").Value = '3.497.18'
$ws.Range("D19").Value = '12.25'
$ws.Range("D20").Value = '0.120'
$ws.Range("D21").Value = '18.17'
$ws.Range("D22").Value = '0.942'
$ws.Range("D23").Value = '18.08'
$ws.Range("D24").Value = '5.26'
$ws.Range("D25").Value = '97.64'
$ws.Range("D26").Value = '4.24'
$ws.Range("D27").Value = '2.85'
$ws.Range("D28").Value = '9.95'
$ws.Range("D29").Value = '9.21'
$ws.Range("D30").Value = '32.21'
$ws.Range("D31").Value = '4.09'
$ws.Range("D32").Value = '6.63'
$ws.Range("D33").Value = '581.80'
$ws.Range("D34").Value = '11.48'
$ws.Range("D35").Value = '60.91'
$ws.Range("D36").Value = '0.108'
$ws.Range("D37").Value = '3.729.40'
$ws.Range("D38").Value = '0.999'
$ws.Range("D39").Value = '3.77'
$ws.Range("D40").Value = '0.0₃0776'
$ws.Range("D41").Value = '3.54'
$ws.Range("D42").Value = '0.370'
$ws.Range("D43").Value = '2.80'
$ws.Range("D44").Value = '0.132'
$ws.Range("D45").Value = '33.63'
$ws.Range("D46").Value = '0.0436'
$ws.Range("D47").Value = '3.30'
$ws.Range("D48").Value = '2.80'
$ws.Range("D49").Value = '0.133'
$ws.Range("D50").Value = '0.999'
$ws.Range("D51").Value = '8.09'

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = '  -2.43%  '
$ws.Range("E3").Value = '  -3.90%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("E5").Value = '  -4.97%  '
$ws.Range("E6").Value = '  -6.12%  '
$ws.Range("E7").Value = '  -3.80%  '
$ws.Range("E8").Value = '  -3.72%  '
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("E10").Value = '  -8.07%  '
$ws.Range("E11").Value = '  -4.75%  '
$ws.Range("E12").Value = '  -5.21%  '
$ws.Range("E13").Value = '  -7.90%  '
$ws.Range("E14").Value = '  -5.95%  '
$ws.Range("E15").Value = '  -3.87%  '
$ws.Range("E16").Value = '  -7.75%  '
$ws.Range("E17").Value = '  -2.82%  '
$ws.Range("E18").Value = '  -3.42%  '
$ws.Range("E19").Value = '  -4.14%  '
$ws.Range("E20").Value = '  -2.25%  '
$ws.Range("E21").Value = '  -4.79%  '
$ws.Range("E22").Value = '  -5.77%  '
$ws.Range("E23").Value = '  -5.38%  '
$ws.Range("E24").Value = '  -3.04%  '
$ws.Range("E25").Value = '  -7.88%  '
$ws.Range("E26").Value = '  -8.15%  '
$ws.Range("E27").Value = '  -6.16%  '
$ws.Range("E28").Value = '  -6.17%  '
$ws.Range("E29").Value = '  -8.57%  '
$ws.Range("E30").Value = '  -5.71%  '
$ws.Range("E31").Value = '  -10.03%  '
$ws.Range("E32").Value = '  -8.79%  '
$ws.Range("E33").Value = '  +8.98%  '
$ws.Range("E34").Value = '  -6.18%  '
$ws.Range("E35").Value = '  -4.18%  '
$ws.Range("E36").Value = '  -6.32%  '
$ws.Range("E37").Value = '  -4.88%  '
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("E39").Value = '  +38.92%  '
$ws.Range("E40").Value = '  -9.99%  '
$ws.Range("E41").Value = '  -2.21%  '
$ws.Range("E42").Value = '  -5.32%  '
$ws.Range("E43").Value = '  -8.25%  '
$ws.Range("E44").Value = '  -5.99%  '
$ws.Range("E45").Value = '  -8.80%  '
$ws.Range("E46").Value = '  -6.57%  '
$ws.Range("E47").Value = '  -5.00%  '
$ws.Range("E48").Value = '  -8.08%  '
$ws.Range("E49").Value = '  -5.37%  '
$ws.Range("E50").Value = '  -0.37%  '
$ws.Range("E51").Value = '  -6.45%  '
